# Re-applies the latest cryptos-list scrape onto the sheet.
# Price/volume cells are stored as *text* in the source data (note several
# price cells use dotted thousands separators like '29.319.73', which aren't
# valid numbers at all) so plain numeric-looking values are entered with a
# leading apostrophe to force text entry, matching the original formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.319.73'
$ws.Range('E2').Value = '  +2.77%  '
$ws.Range('D3').Value = '1.895.08'
$ws.Range('E3').Value = '  +0.88%  '
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').Value = '''314.62'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('D7').Value = '''0.5148'
$ws.Range('E7').Value = '  +0.77%  '
$ws.Range('D8').Value = '''0.3925'
$ws.Range('E8').Value = '  -0.36%  '
$ws.Range('D9').Value = '''0.08430'
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('D10').Value = '''42.39'
$ws.Range('E10').Value = '  +1.85%  '
$ws.Range('D11').Value = '''1.116'
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('D12').Value = '''6.265'
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('D13').Value = '1.897.50'
$ws.Range('E13').Value = '  +1.36%  '
$ws.Range('D14').Value = '''20.70'
$ws.Range('D15').Value = '''7.293'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').Value = '''1.004'
$ws.Range('E16').Value = '  -0.40%  '
$ws.Range('D17').Value = '''93.22'
$ws.Range('E17').Value = '  +2.36%  '
$ws.Range('D18').Value = '''0.00001104'
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').Value = '''0.06731'
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('D20').Value = '''17.85'
$ws.Range('E20').Value = '  +0.76%  '
$ws.Range('D21').Value = '''1.002'
$ws.Range('E21').Value = '  -0.39%  '
$ws.Range('D22').Value = '''6.014'
$ws.Range('E22').Value = '  +1.02%  '
$ws.Range('D23').Value = '29.340.35'
$ws.Range('E23').Value = '  +2.79%  '
$ws.Range('D24').Value = '''11.15'
$ws.Range('E24').Value = '  +0.32%  '
$ws.Range('D25').Value = '''2.215'
$ws.Range('E25').Value = '  -1.88%  '
$ws.Range('D26').Value = '2.111.95'
$ws.Range('E26').Value = '  +1.17%  '
$ws.Range('D27').Value = '''159.29'
$ws.Range('E27').Value = '  -1.14%  '
$ws.Range('D28').Value = '''20.87'
$ws.Range('E28').Value = '  +0.79%  '
$ws.Range('D29').Value = '''2.433'
$ws.Range('E29').Value = '  +2.17%  '
$ws.Range('D30').Value = '''127.37'
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '''1.061'
$ws.Range('E31').Value = '  +1.15%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '''0.1050'
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').Value = '''6.119'
$ws.Range('E33').Value = '  +5.80%  '
$ws.Range('D34').Value = '''3.658'
$ws.Range('E34').Value = '  +1.41%  '
$ws.Range('D35').Value = '''0.02481'
$ws.Range('E35').Value = '  +1.29%  '
$ws.Range('D36').Value = '''0.06577'
$ws.Range('E36').Value = '  +0.99%  '
$ws.Range('D37').Value = '''0.2199'
$ws.Range('E37').Value = '  +0.36%  '
$ws.Range('D38').Value = '''9.049'
$ws.Range('E38').Value = '  +1.35%  '
$ws.Range('D39').Value = '''5.193'
$ws.Range('E39').Value = '  +2.30%  '
$ws.Range('D40').Value = '''1.229'
$ws.Range('E40').Value = '  +2.81%  '
$ws.Range('D41').Value = '''0.6520'
$ws.Range('E41').Value = '  +1.24%  '
$ws.Range('D42').Value = '''1.234'
$ws.Range('E42').Value = '  -2.03%  '
$ws.Range('D43').Value = '''11.28'
$ws.Range('E43').Value = '  +1.02%  '
$ws.Range('D44').Value = '''0.6072'
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').Value = '''13.18'
$ws.Range('E45').Value = '  +1.36%  '
$ws.Range('D46').Value = '''3.672'
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('D47').Value = '''2.051'
$ws.Range('E47').Value = '  +2.15%  '
$ws.Range('D48').Value = '''1.228'
$ws.Range('E48').Value = '  +1.60%  '
$ws.Range('D49').Value = '''123.46'
$ws.Range('E50').Value = '  -3.53%  '
$ws.Range('D51').Value = '''77.71'
$ws.Range('E51').Value = '  +0.78%  '
